# Append the July 28, 2020 SSA data row (row 59) to the "out_vars" sheet,
# matching the existing table's layout (Fecha, Confirmados, Negativos,
# Sospechosos, Defunciones, Porcentaje hospitalizados).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("out_vars")

$row = 59

# Numeric columns first.
$ws.Cells.Item($row, 2).Value = 402697
$ws.Cells.Item($row, 3).Value = 449854
$ws.Cells.Item($row, 4).Value = 87538
$ws.Cells.Item($row, 5).Value = 44876
$ws.Cells.Item($row, 6).Value = 27.63

# The date column is stored as plain text in this sheet (e.g. "2020-07-27"),
# not as an Excel date serial, so force text formatting before assigning the
# value and then clear the formatting again so the cell keeps the sheet's
# default (unstyled) look, same as every other row in the column.
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2020-07-28"
$dateCell.ClearFormats()
